$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 21:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 831086
$ws.Range("C4").Value = 12342
$ws.Range("D4").Value = 83462
$ws.Range("E4").Value = 701375
$ws.Range("F4").Value = 14014
$ws.Range("G4").Value = 931
$ws.Range("H4").Value = 46249

# Row 14 - Brasil
$ws.Range("B14").Value = 45757
$ws.Range("C14").Value = 2678
$ws.Range("E14").Value = 18526
$ws.Range("G14").Value = 165
$ws.Range("H14").Value = 2906

# Row 21 - Peru
$ws.Range("B21").Value = 19250
$ws.Range("C21").Value = 1413
$ws.Range("D21").Value = 7027
$ws.Range("E21").Value = 11693
$ws.Range("F21").Value = 396
$ws.Range("G21").Value = 46
$ws.Range("H21").Value = 530

# Row 96 - Costa Rica
$ws.Range("B96").Value = 681
$ws.Range("C96").Value = 12
$ws.Range("D96").Value = 180
$ws.Range("E96").Value = 495
